$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 97
$ws.Range("H97").Value = 1889.9
$ws.Range("J97").Value = 1889.9
$ws.Range("L97").Value = 5669.700000000001
$ws.Range("N97").Value = -6661.700000000001

# Row 100
$ws.Range("H100").Value = 2389.0527
$ws.Range("I100").Value = 1273.25
$ws.Range("K100").Value = 1273.25
$ws.Range("M100").Value = -732.25

# Row 103
$ws.Range("H103").Value = 50000670
$ws.Range("I103").Value = 754.125
$ws.Range("J103").Value = 250000350
$ws.Range("K103").Value = 2262.375
$ws.Range("L103").Value = 750001050
$ws.Range("M103").Value = -1676.375
$ws.Range("N103").Value = -750002222

# Row 106
$ws.Range("H106").Value = 3649.3333
$ws.Range("I106").Value = 3649.3333
$ws.Range("K106").Value = 3649.3333
$ws.Range("M106").Value = -3018.3333

# Row 110
$ws.Range("H110").Value = 37950
$ws.Range("J110").Value = 37950
$ws.Range("L110").Value = 37950
$ws.Range("N110").Value = -46130

# Row 112
$ws.Range("H112").Value = 2800.6
$ws.Range("J112").Value = 2949.4055
$ws.Range("L112").Value = 8848.216499999999
$ws.Range("N112").Value = -11064.2165

# Row 113
$ws.Range("H113").Value = 3036.9473
$ws.Range("I113").Value = 2483.3333
$ws.Range("J113").Value = 3292.4614
$ws.Range("K113").Value = 2483.3333
$ws.Range("L113").Value = 3292.4614
$ws.Range("M113").Value = 770.6667000000002
$ws.Range("N113").Value = -9800.4614

# Row 115
$ws.Range("H115").Value = 957.2778
$ws.Range("I115").Value = 611.6667
$ws.Range("J115").Value = 1130.0834
$ws.Range("K115").Value = 1835.0001
$ws.Range("L115").Value = 3390.2502
$ws.Range("M115").Value = -268.0001
$ws.Range("N115").Value = -6524.2502

# Row 118
$ws.Range("H118").Value = 922.2857
$ws.Range("I118").Value = 423
$ws.Range("J118").Value = 2520
$ws.Range("K118").Value = 1269
$ws.Range("L118").Value = 7560
$ws.Range("M118").Value = 388
$ws.Range("N118").Value = -10874

# Row 141
$ws.Range("H141").Value = 5553.846
$ws.Range("I141").Value = 5917.727
$ws.Range("J141").Value = 3552.5
$ws.Range("K141").Value = 17753.181
$ws.Range("L141").Value = 10657.5
$ws.Range("M141").Value = -12573.181
$ws.Range("N141").Value = -21017.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1561.7059
$ws.Range("I61").Value = 1494.9756
$ws.Range("J61").Value = 1835.3
$ws.Range("K61").Value = 1494.9756
$ws.Range("L61").Value = 1835.3
$ws.Range("M61").Value = -1282.9756
$ws.Range("N61").Value = -2259.3

# Row 102
$ws.Range("H102").Value = 5765.0713
$ws.Range("I102").Value = 5390
$ws.Range("J102").Value = 6702.75
$ws.Range("K102").Value = 5390
$ws.Range("L102").Value = 6702.75
$ws.Range("M102").Value = -3768
$ws.Range("N102").Value = -9946.75

# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# Row 132
$ws.Range("H132").Value = 1801.72
$ws.Range("I132").Value = 1291.2858
$ws.Range("J132").Value = 4481.5
$ws.Range("K132").Value = 3873.8574
$ws.Range("L132").Value = 13444.5
$ws.Range("M132").Value = -1343.8574
$ws.Range("N132").Value = -18504.5

# Row 136
$ws.Range("H136").Value = 1561.7059
$ws.Range("I136").Value = 1494.9756
$ws.Range("J136").Value = 1835.3
$ws.Range("K136").Value = 4484.9268
$ws.Range("L136").Value = 5505.9
$ws.Range("M136").Value = -1934.9268
$ws.Range("N136").Value = -10605.9

# Row 138
$ws.Range("H138").Value = 47011.43
$ws.Range("J138").Value = 47011.43
$ws.Range("L138").Value = 47011.43
$ws.Range("N138").Value = -57291.43

# Row 141
$ws.Range("H141").Value = 66000
$ws.Range("J141").Value = 66000
$ws.Range("L141").Value = 66000
$ws.Range("N141").Value = -76360

$ws = $wb.Worksheets.Item("BSM")
# Row 102
$ws.Range("H102").Value = 14852
$ws.Range("I102").Value = 14852
$ws.Range("K102").Value = 14852
$ws.Range("M102").Value = -11607

# Row 103
$ws.Range("H103").Value = 25828.5
$ws.Range("J103").Value = 25828.5
$ws.Range("L103").Value = 25828.5
$ws.Range("N103").Value = -28172.5

# Row 105
$ws.Range("H105").Value = 1659.375
$ws.Range("I105").Value = 1590
$ws.Range("J105").Value = 1690.909
$ws.Range("K105").Value = 1590
$ws.Range("L105").Value = 1690.909
$ws.Range("M105").Value = 157
$ws.Range("N105").Value = -5184.909

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# Row 140
$ws.Range("H140").Value = 76480
$ws.Range("J140").Value = 76480
$ws.Range("L140").Value = 76480
$ws.Range("N140").Value = -86840

# Row 141
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 43
$ws.Range("H43").Value = 27655.834
$ws.Range("J43").Value = 27655.834
$ws.Range("L43").Value = 27655.834
$ws.Range("N43").Value = -28023.834

# Row 68
$ws.Range("H68").Value = 17701
$ws.Range("J68").Value = 17701
$ws.Range("L68").Value = 17701
$ws.Range("N68").Value = -19199

# Row 71
$ws.Range("H71").Value = 17701
$ws.Range("J71").Value = 17701
$ws.Range("L71").Value = 53103
$ws.Range("N71").Value = -60591

# Row 101
$ws.Range("H101").Value = 27655.834
$ws.Range("J101").Value = 27655.834
$ws.Range("L101").Value = 27655.834
$ws.Range("N101").Value = -34145.834

# Row 102
$ws.Range("H102").Value = 48500
$ws.Range("J102").Value = 48500
$ws.Range("L102").Value = 48500
$ws.Range("N102").Value = -53368

# Row 103
$ws.Range("H103").Value = 2804.8
$ws.Range("I103").Value = 2804.8
$ws.Range("K103").Value = 2804.8
$ws.Range("M103").Value = -1632.8

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 30794.416
$ws.Range("I132").Value = 42720.36
$ws.Range("J132").Value = 3690
$ws.Range("K132").Value = 128161.08
$ws.Range("L132").Value = 11070
$ws.Range("M132").Value = -125631.08
$ws.Range("N132").Value = -16130

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

# Row 58
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 10000
$ws.Range("M58").Value = -9740

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Row 100
$ws.Range("H100").Value = 2974.1365
$ws.Range("I100").Value = 2190.6
$ws.Range("J100").Value = 3627.0833
$ws.Range("K100").Value = 2190.6
$ws.Range("L100").Value = 3627.0833
$ws.Range("M100").Value = -1649.6
$ws.Range("N100").Value = -4709.0833

# Row 136
$ws.Range("H136").Value = 4765
$ws.Range("I136").Value = 2579.2727
$ws.Range("J136").Value = 8199.714
$ws.Range("K136").Value = 7737.8181
$ws.Range("L136").Value = 24599.142
$ws.Range("M136").Value = -5187.8181
$ws.Range("N136").Value = -29699.142

$ws = $wb.Worksheets.Item("WVR")
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 107
$ws.Range("H107").Value = 832.61536
$ws.Range("I107").Value = 852.4
$ws.Range("J107").Value = 766.6667
$ws.Range("K107").Value = 2557.2
$ws.Range("L107").Value = 2300.0001
$ws.Range("M107").Value = -637.1999999999998
$ws.Range("N107").Value = -6140.0001

# Row 113
$ws.Range("H113").Value = 32258776
$ws.Range("I113").Value = 551.26086
$ws.Range("J113").Value = 125001176
$ws.Range("K113").Value = 1653.78258
$ws.Range("L113").Value = 375003528
$ws.Range("M113").Value = 516.2174199999999
$ws.Range("N113").Value = -375007868
